$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: I1 = "I0", J1 = "IF" (styled like the other header cells, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-26
$values = @{
    2  = @(10, 10)
    3  = @(9, 9)
    4  = @(7, 8)
    5  = @(9, 9)
    6  = @(8, 8)
    7  = @(6, 6)
    8  = @(6, 7)
    9  = @(8, 8)
    10 = @(8, 9)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(9, 9)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(9, 9)
    20 = @(2, 2)
    21 = @(8, 8)
    22 = @(6, 6)
    23 = @(5, 5)
    24 = @(5, 5)
    25 = @(5, 5)
    26 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

$wb.Save()
